{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// The targets are the (now trailing, no-longer-needed) empty paragraph and\n// the two footer paragraphs that directly follow the final bibliography\n// entry (\"... Autor: Bertero, C. O. Editora: ATLAS\").\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet bertoroIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais Autor: Bertero\") !== -1) {\n    bertoroIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (bertoroIndex !== -1 && items[bertoroIndex + 1] && items[bertoroIndex + 1].text === \"\") {\n  // the blank paragraph right after the Bertero bibliography line\n  toDelete.push(items[bertoroIndex + 1]);\n}\nfor (const para of items) {\n  if (targetTexts.indexOf(para.text) !== -1) {\n    toDelete.push(para);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$bertero = \"Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais Autor: Bertero, C. O. Editora: ATLAS\"\n$jupiter = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyright = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Snapshot paragraphs (live Range objects) and their trimmed text up front,\n# since collection indices shift as paragraphs are removed.\n$paras = @($d.Paragraphs)\n$texts = @()\nforeach ($p in $paras) {\n    $texts += ,$p.Range.Text.TrimEnd([char]0x0D, [char]0x07)\n}\n\n$berteroIndex = -1\nfor ($i = 0; $i -lt $texts.Count; $i++) {\n    if ($texts[$i] -eq $bertero) {\n        $berteroIndex = $i\n        break\n    }\n}\n\n# Collect the indices to remove: the blank paragraph right after the\n# \"Bertero\" bibliography entry, plus the \"Ver no Jupiter...\" and the\n# \"\u00a9 2020...\" footer paragraphs.\n$toRemove = New-Object System.Collections.Generic.List[int]\nif ($berteroIndex -ge 0 -and ($berteroIndex + 1) -lt $texts.Count -and $texts[$berteroIndex + 1] -eq \"\") {\n    [void]$toRemove.Add($berteroIndex + 1)\n}\nfor ($i = 0; $i -lt $texts.Count; $i++) {\n    if ($texts[$i] -eq $jupiter -or $texts[$i] -eq $copyright) {\n        [void]$toRemove.Add($i)\n    }\n}\n\n# Delete from the bottom up so earlier indices stay valid.\n$sorted = $toRemove | Sort-Object -Descending -Unique\nforeach ($idx in $sorted) {\n    $paras[$idx].Range.Delete()\n}\n"}
